# Fruta / hortaliza, semanal
# Weekly data refresh for Granada (Vega Central Mapocho de Santiago) sheet.
# Rewrites data rows 2-22 (adds rows 21-22 for the new week) and updates
# the used range dimension accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r2 = @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44305, 13, 'Fruta', 100104, 'Frutos de pepita', 100104001, 'Granada', 'Wonderfull', 'Primera', 50, 18000, 18000, 18000, '$/caja 15 kilos granel', 'Región de O''Higgins', 1200, 15)
$r3 = @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44305, 13, 'Fruta', 100104, 'Frutos de pepita', 100104001, 'Granada', 'Wonderfull', 'Segunda', 60, 15000, 15000, 15000, '$/caja 15 kilos granel', 'Región de O''Higgins', 1000, 15)
$r4 = @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44678, 13, 'Fruta', 100104, 'Frutos de pepita', 100104001, 'Granada', 'Sin especificar', 'Especial', 290, 15000, 15000, 15000, '$/caja 15 kilos granel', 'Región de O''Higgins', 1000, 15)
$r5 = @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44678, 13, 'Fruta', 100104, 'Frutos de pepita', 100104001, 'Granada', 'Sin especificar', 'Primera', 220, 12000, 12000, 12000, '$/caja 15 kilos granel', 'Región de O''Higgins', 800, 15)
$r6 = @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44309, 13, 'Fruta', 100104, 'Frutos de pepita', 100104001, 'Granada', 'Wonderfull', 'Primera', 40, 18000, 18000, 18000, '$/caja 15 kilos granel', 'Región de O''Higgins', 1200, 15)
$r7 = @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44309, 13, 'Fruta', 100104, 'Frutos de pepita', 100104001, 'Granada', 'Wonderfull', 'Segunda', 70, 15000, 15000, 15000, '$/caja 15 kilos granel', 'Región de O''Higgins', 1000, 15)
$r8 = @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44658, 13, 'Fruta', 100104, 'Frutos de pepita', 100104001, 'Granada', 'Sin especificar', 'Especial', 280, 21600, 21600, 21600, '$/caja 18 kilos granel', 'Provincia de Limarí', 1200, 18)
$r9 = @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44658, 13, 'Fruta', 100104, 'Frutos de pepita', 100104001, 'Granada', 'Sin especificar', 'Primera', 330, 16200, 16200, 16200, '$/caja 18 kilos granel', 'Provincia de Limarí', 900, 18)
$r10 = @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44658, 13, 'Fruta', 100104, 'Frutos de pepita', 100104001, 'Granada', 'Sin especificar', 'Segunda', 220, 14400, 14400, 14400, '$/caja 18 kilos granel', 'Provincia de Limarí', 800, 18)
$r11 = @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44649, 13, 'Fruta', 100104, 'Frutos de pepita', 100104001, 'Granada', 'Sin especificar', 'Especial', 220, 21600, 21600, 21600, '$/caja 18 kilos granel', 'Provincia de Limarí', 1200, 18)
$r12 = @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44649, 13, 'Fruta', 100104, 'Frutos de pepita', 100104001, 'Granada', 'Sin especificar', 'Primera', 250, 16200, 16200, 16200, '$/caja 18 kilos granel', 'Provincia de Limarí', 900, 18)
$r13 = @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44649, 13, 'Fruta', 100104, 'Frutos de pepita', 100104001, 'Granada', 'Sin especificar', 'Segunda', 180, 14400, 14400, 14400, '$/caja 18 kilos granel', 'Provincia de Limarí', 800, 18)
$r14 = @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44664, 13, 'Fruta', 100104, 'Frutos de pepita', 100104001, 'Granada', 'Sin especificar', 'Especial', 300, 21600, 21600, 21600, '$/caja 18 kilos granel', 'Provincia de Limarí', 1200, 18)
$r15 = @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44664, 13, 'Fruta', 100104, 'Frutos de pepita', 100104001, 'Granada', 'Sin especificar', 'Primera', 250, 18000, 18000, 18000, '$/caja 18 kilos granel', 'Provincia de Limarí', 1000, 18)
$r16 = @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44664, 13, 'Fruta', 100104, 'Frutos de pepita', 100104001, 'Granada', 'Sin especificar', 'Segunda', 250, 16000, 16000, 16000, '$/caja 18 kilos granel', 'Provincia de Limarí', 889, 18)
$r17 = @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44285, 13, 'Fruta', 100104, 'Frutos de pepita', 100104001, 'Granada', 'Wonderfull', 'Especial', 40, 18000, 18000, 18000, '$/caja 15 kilos empedrada', 'Provincia del Elquí', 1200, 15)
$r18 = @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44285, 13, 'Fruta', 100104, 'Frutos de pepita', 100104001, 'Granada', 'Wonderfull', 'Primera', 90, 15000, 15000, 15000, '$/caja 15 kilos empedrada', 'Provincia del Elquí', 1000, 15)
$r19 = @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44285, 13, 'Fruta', 100104, 'Frutos de pepita', 100104001, 'Granada', 'Wonderfull', 'Segunda', 75, 12000, 12000, 12000, '$/caja 15 kilos empedrada', 'Provincia del Elquí', 800, 15)
$r20 = @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44644, 13, 'Fruta', 100104, 'Frutos de pepita', 100104001, 'Granada', 'Sin especificar', 'Especial', 180, 18000, 18000, 18000, '$/caja 15 kilos granel', 'Provincia de Limarí', 1200, 15)
$r21 = @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44644, 13, 'Fruta', 100104, 'Frutos de pepita', 100104001, 'Granada', 'Sin especificar', 'Primera', 220, 13500, 13500, 13500, '$/caja 15 kilos granel', 'Provincia de Limarí', 900, 15)
$r22 = @(9, 'Vega Central Mapocho de Santiago', 'Metropolitana', 44644, 13, 'Fruta', 100104, 'Frutos de pepita', 100104001, 'Granada', 'Sin especificar', 'Segunda', 290, 12000, 12000, 12000, '$/caja 15 kilos granel', 'Provincia de Limarí', 800, 15)

$rows = @($r2, $r3, $r4, $r5, $r6, $r7, $r8, $r9, $r10, $r11, $r12, $r13, $r14, $r15, $r16, $r17, $r18, $r19, $r20, $r21, $r22)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $rowData = $rows[$i]
    $excelRow = $startRow + $i
    for ($c = 1; $c -le $rowData.Count; $c++) {
        $ws.Cells.Item($excelRow, $c).Value = $rowData[$c - 1]
    }
    # Column D (Fecha) keeps the workbook's date/time display format.
    $ws.Cells.Item($excelRow, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
